# İş Takip Güncellemesi - 24.07.2025 14:54:32
# Populate the "İş Takip Listesi" sheet (the first / active sheet) with a
# header row plus two data rows, exactly as in the source update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the whole target range to Text format first so values that look
# like numbers ("5", "120") or dates ("2025-01-29") are written as literal
# text strings instead of being auto-converted to numbers/dates.
$ws.Range("A1:L3").NumberFormat = "@"

$headers = @(
    "Tamamlandı",
    "İL",
    "İLÇE",
    "BİRİM",
    "UYGULAMA",
    "GÖREVLİ PERSONELLER",
    "PARSEL SAYISI",
    "ALAN(Ha)",
    "İHALELİ/MÜDÜRLÜK",
    "İŞE BAŞLAMA/YER TESLİMİ",
    "İHALE BİTİŞ TARİHİ",
    "DURUMU"
)

$row2 = @(
    "HAYIR",
    "Adana",
    "Akdeniz",
    "ghg",
    "GÜNCELLEME",
    "ESMEN TOKALI (K.Mühendisi), SERVET ATA (K.Mühendisi)",
    "",
    "5",
    "120",
    "İhaleli",
    "2025-01-29",
    "2025-07-24"
)

$row3 = @(
    "HAYIR",
    "Adana",
    "Erdemli",
    "hjhj",
    "2/B",
    "MUSTAFA DIKI (S. Mühendis), HAKAN ÖZEL (K.Teknisyeni)",
    "",
    "15",
    "155",
    "İhaleli",
    "2025-02-01",
    "2025-07-24"
)

$data = @($headers, $row2, $row3)

for ($r = 1; $r -le $data.Length; $r++) {
    $rowValues = $data[$r - 1]
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
